# Update "南宁-漫展信息" workbook: remove the expired first event (南宁·第一届异次元动漫嘉年华),
# shift the remaining events up by one row, renumber the index column, and refresh the
# "想去人数" (interest count) figures that changed since the last scrape.
# This applies identically to the "展览" (sheet 1) and "全部类型" (sheet 4) worksheets,
# which carry duplicate data in this workbook.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Remove the whole first data row (row 2); this shifts rows 3-6 up to 2-5
    # and keeps all formatting/styles intact, just like a manual row delete.
    $ws.Rows.Item(2).Delete()

    # Renumber the index column (A) for the four remaining records.
    $ws.Cells.Item(2, 1).Value = 1
    $ws.Cells.Item(3, 1).Value = 2
    $ws.Cells.Item(4, 1).Value = 3
    $ws.Cells.Item(5, 1).Value = 4

    # Refresh the "想去人数" (interest count) numbers that changed in this scrape.
    $ws.Cells.Item(2, 6).Value = 1688
    $ws.Cells.Item(3, 6).Value = 7793
    $ws.Cells.Item(5, 6).Value = 233
}
